# Update "想去人数" (interest count) figures in the 展览 and 全部类型 sheets
# to reflect the newly generated numbers from the commit.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1): column F holds 想去人数
$wsExhibit.Range("F3").Value  = 61
$wsExhibit.Range("F4").Value  = 1480
$wsExhibit.Range("F5").Value  = 567
$wsExhibit.Range("F7").Value  = 10979
$wsExhibit.Range("F8").Value  = 10979
$wsExhibit.Range("F11").Value = 312
$wsExhibit.Range("F13").Value = 748
$wsExhibit.Range("F15").Value = 12704

# 全部类型 sheet (sheet4): same events, shifted down by one row
$wsAll.Range("F4").Value  = 61
$wsAll.Range("F5").Value  = 1480
$wsAll.Range("F6").Value  = 567
$wsAll.Range("F8").Value  = 10979
$wsAll.Range("F9").Value  = 10979
$wsAll.Range("F12").Value = 312
$wsAll.Range("F14").Value = 748
$wsAll.Range("F16").Value = 12704
